# Generate Report for Handback
# Adds a new "version mismatch" handback result for the
# 6a65ad79-ad47-4002-acab-b46267cb6f93 item on both the zh-cn and de-de
# status sheets: a hyperlinked "Latest Target File" entry, the matching
# target xlf name, a new "Latest Handback DateTime" stamp, and an
# "Error Detail" message - plus widening the Error Detail column.

$wb = $excel.ActiveWorkbook

$latestHandbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d3bd5495f3614cdf58975d2d5b93e314c73b0041/e2e/6a65ad79-ad47-4002-acab-b46267cb6f93.md"
$handbackMdName = "6a65ad79-ad47-4002-acab-b46267cb6f93.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bd2d852fe3934c0fdee831093c5c9acd909e1d9f/e2e/6a65ad79-ad47-4002-acab-b46267cb6f93.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d3bd5495f3614cdf58975d2d5b93e314c73b0041/e2e/6a65ad79-ad47-4002-acab-b46267cb6f93.md."

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I8"), $latestHandbackUrl, "", "", $handbackMdName)
$wsZhCn.Range("J8").Value = "6a65ad79-ad47-4002-acab-b46267cb6f93.96e5fcf36d338a61e40326bf6cf0647ba5e8c465.zh-cn.xlf"
$wsZhCn.Range("K8").Value = "2016-08-21 02:47:40"
$wsZhCn.Range("P8").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I8"), $latestHandbackUrl, "", "", $handbackMdName)
$wsDeDe.Range("J8").Value = "6a65ad79-ad47-4002-acab-b46267cb6f93.96e5fcf36d338a61e40326bf6cf0647ba5e8c465.de-de.xlf"
$wsDeDe.Range("K8").Value = "2016-08-21 02:47:46"
$wsDeDe.Range("P8").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
